# commit pyscript update 4
# 1. add timeout column/param handling isn't represented in the sheet grid
# 2. add token handling -> new "savetoken" column (M) on Sheet1, plus its two
#    per-row values, and repoint the "tokenname" sample values to the new
#    "bdh" env-var style token name
# 3. environment variable support -> the hard-coded host strings in row 3 of
#    Sheet1 are replaced by the short "bdh" placeholder (env-var driven)
# 4. excel data handled as json -> no grid-visible change beyond the above

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: update existing values --------------------------------------
# tokenname sample on row 2 now references the "bdh" env/token name instead
# of the literal header name
$ws1.Range("K2").Value = "bdh"
# url sample on row 3 collapses to the "bdh" placeholder (env var) instead of
# the hard-coded " dev01.bdhlan.com:8080/bdhsystem" host string
$ws1.Range("C3").Value = "bdh"

# --- Sheet1: new "savetoken" column (M) -----------------------------------
$ws1.Range("M1").Value = "savetoken"
$ws1.Range("M2").Value = "header_data_token"
$ws1.Range("M3").Value = "param_location_datas_data_token"

# --- Sheet2: printer/page setup now explicitly configured ----------------
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- restore the on-screen selections the saved file records -------------
$ws2.Activate() | Out-Null
$ws2.Range("G2").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("M1:M3").Select() | Out-Null
